# Atualizacao de bases das ligas, do dia: 19-06-2024 as 21:51
#
# This update re-shuffles the per-match data (everything except the
# row's rank in column A and the round's Date in column D) among a
# number of rows on the single worksheet. For every affected row the
# columns B..AD (match id, league, home/away team, scores, odds, etc.)
# are moved to a different row while column A (rank) and column D
# (date) stay where they are.
#
# The affected rows form the following permutation cycles (row <- row
# that "donates" its B:AD data to it):
#   19 <- 20, 20 <- 21, 21 <- 22, 22 <- 19        (4-cycle)
#   23 <- 24, 24 <- 25, 25 <- 23                  (3-cycle)
#   36 <-> 37
#   38 <-> 39
#   40 <-> 41
#   43 <-> 44
#   47 <-> 48
#   108 <-> 109
#   210 <-> 211

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B$row`:AD$row").Value2
}

function Set-RowData($row, $data) {
    $ws.Range("B$row`:AD$row").Value2 = $data
}

# ---- Cyclic permutations (length > 2) ----
$cycles = @(
    @(19, 20, 21, 22),
    @(23, 24, 25)
)

foreach ($cycle in $cycles) {
    # snapshot all rows in the cycle before writing any of them
    $snapshot = @{}
    foreach ($r in $cycle) {
        $snapshot[$r] = Get-RowData $r
    }
    # new[cycle[i]] = old[cycle[i+1]] (wrapping around)
    $n = $cycle.Length
    for ($i = 0; $i -lt $n; $i++) {
        $dst = $cycle[$i]
        $src = $cycle[($i + 1) % $n]
        Set-RowData $dst $snapshot[$src]
    }
}

# ---- Simple pairwise swaps ----
$pairs = @(
    @(36, 37),
    @(38, 39),
    @(40, 41),
    @(43, 44),
    @(47, 48),
    @(108, 109),
    @(210, 211)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $d1 = Get-RowData $r1
    $d2 = Get-RowData $r2
    Set-RowData $r1 $d2
    Set-RowData $r2 $d1
}
